# Update for NFL 2020 season
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "nfl_tm": add the Las Vegas Raiders franchise row and mark the
# old Oakland Raiders record as inactive (status 1 -> 0).
# ---------------------------------------------------------------------
$tm = $wb.Worksheets.Item("nfl_tm")

# Oakland Raiders (row 26) is no longer the active Raiders record.
$tm.Cells.Item(26, 7).Value = 0

# Add the new Las Vegas Raiders row (row 37).
$tm.Cells.Item(37, 1).Value = 36
$tm.Cells.Item(37, 2).Value = "LVR"
$tm.Cells.Item(37, 3).Value = "Las Vegas Raiders"
$tm.Cells.Item(37, 4).Value = "Las Vegas"
$tm.Cells.Item(37, 5).Value = "LAV"
$tm.Cells.Item(37, 6).Value = "Las Vegas Raiders"
$tm.Cells.Item(37, 7).Value = 1

# Column C (tm_name_full) is widened to fit its contents.
$tm.Columns.Item(3).AutoFit() | Out-Null

# ---------------------------------------------------------------------
# Sheet "nfl_game_result": append the week 1, 2020 season schedule.
# ---------------------------------------------------------------------
$gr = $wb.Worksheets.Item("nfl_game_result")

$newGames = @(
    @(2116, 1, 2020, 1, 1, "HOU", "KCY"),
    @(2117, 1, 2020, 1, 1, "SEA", "ATL"),
    @(2118, 1, 2020, 1, 1, "NYJ", "BUF"),
    @(2119, 1, 2020, 1, 1, "CHI", "DET"),
    @(2120, 1, 2020, 1, 1, "CLE", "BAL"),
    @(2121, 1, 2020, 1, 1, "GBY", "MIN"),
    @(2122, 1, 2020, 1, 1, "IND", "JAX"),
    @(2123, 1, 2020, 1, 1, "LAV", "CAR"),
    @(2124, 1, 2020, 1, 1, "MIA", "NWE"),
    @(2125, 1, 2020, 1, 1, "PHI", "WAS"),
    @(2126, 1, 2020, 1, 1, "LAC", "CIN"),
    @(2127, 1, 2020, 1, 1, "TBY", "NOR"),
    @(2128, 1, 2020, 1, 1, "ARI", "SFO"),
    @(2129, 1, 2020, 1, 1, "DAL", "LAR"),
    @(2130, 1, 2020, 1, 1, "PIT", "NYG"),
    @(2131, 1, 2020, 1, 1, "TEN", "DEN")
)

$row = 258
foreach ($g in $newGames) {
    $gr.Cells.Item($row, 1).Value = $g[0]
    $gr.Cells.Item($row, 2).Value = $g[1]
    $gr.Cells.Item($row, 3).Value = $g[2]
    $gr.Cells.Item($row, 4).Value = $g[3]
    $gr.Cells.Item($row, 5).Value = $g[4]
    $gr.Cells.Item($row, 6).Value = $g[5]
    $gr.Cells.Item($row, 7).Value = $g[6]
    $row = $row + 1
}

# Restore the view focus on the sheet / cell that was being edited.
$gr.Activate() | Out-Null
$gr.Range("F266").Select() | Out-Null
